# Commit: "Fruta / hortaliza, semanal"
# A new weekly price-record row is inserted into the data table at sheet
# row 138 (the row right after the header), pushing every following row
# down by one and extending the sheet's data range from A1:R244 to
# A1:R245.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 138

# Insert a brand-new row above the current row 138; all rows from 138
# downward shift down by one position (old 138 -> new 139, ..., old 244 -> new 245).
$ws.Rows($newRow).Insert()

# Fill in the data for the newly inserted row.
$ws.Cells.Item($newRow, 1).Value  = 5
$ws.Cells.Item($newRow, 2).Value  = 'Macroferia Regional de Talca'
$ws.Cells.Item($newRow, 3).Value  = 'Maule'
$ws.Cells.Item($newRow, 4).Value  = 44673
$ws.Cells.Item($newRow, 5).Value  = 7
$ws.Cells.Item($newRow, 6).Value  = 100112021
$ws.Cells.Item($newRow, 7).Value  = 'Ají'
$ws.Cells.Item($newRow, 8).Value  = 'Cacho cabra rojo'
$ws.Cells.Item($newRow, 9).Value  = 'Primera'
$ws.Cells.Item($newRow, 10).Value = 100
$ws.Cells.Item($newRow, 11).Value = 20000
$ws.Cells.Item($newRow, 12).Value = 20000
$ws.Cells.Item($newRow, 13).Value = 20000
$ws.Cells.Item($newRow, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item($newRow, 15).Value = 'Región del Maule'
$ws.Cells.Item($newRow, 16).Value = 800
$ws.Cells.Item($newRow, 17).Value = 25
$ws.Cells.Item($newRow, 18).Value = 'Hortaliza'
